$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 30-35 on the "Artfynd" sheet hold per-observation species records.
# The underlying source export re-synced, so the Id/Taxonsorteringsordning/
# TaxonId/Artnamn/Vetenskapligt namn/Auktor/Antal/Enhet/Alder-Stadium/Ost/Nord
# columns (A,B,E,F,G,H,I,J,K,Q,R) need their refreshed values written in.
# "Antal" (I) is stored as text in the source feed even when numeric-looking,
# so number-format the cell as Text first to avoid Excel auto-converting it.

# Row 30
$ws.Cells.Item(30, 1).Value = 111564898    # A30 Id
$ws.Cells.Item(30, 2).Value = 95532    # B30 Taxonsorteringsordning
$ws.Cells.Item(30, 5).Value = 221945    # E30 TaxonId
$ws.Cells.Item(30, 6).Value = "Revlummer"    # F30 Artnamn
$ws.Cells.Item(30, 7).Value = "Lycopodium annotinum"    # G30 Vetenskapligt namn
$ws.Cells.Item(30, 8).Value = "L."    # H30 Auktor
$ws.Cells.Item(30, 9).Value = ""    # I30 Antal
$ws.Cells.Item(30, 10).Value = ""    # J30 Enhet
$ws.Cells.Item(30, 11).Value = "fullt utvecklade blad"    # K30 Alder-Stadium
$ws.Cells.Item(30, 13).ClearContents()    # M30 Aktivitet (no longer set)
$ws.Cells.Item(30, 17).Value = 561231.5123860433    # Q30 Ost
$ws.Cells.Item(30, 18).Value = 6622624.185033713    # R30 Nord

# Row 31
$ws.Cells.Item(31, 1).Value = 111564905    # A31 Id
$ws.Cells.Item(31, 2).Value = 57578    # B31 Taxonsorteringsordning
$ws.Cells.Item(31, 5).Value = 208250    # E31 TaxonId
$ws.Cells.Item(31, 6).Value = "Åkergroda"    # F31 Artnamn
$ws.Cells.Item(31, 7).Value = "Rana arvalis"    # G31 Vetenskapligt namn
$ws.Cells.Item(31, 8).Value = "Nilsson, 1842"    # H31 Auktor
$ws.Cells.Item(31, 9).NumberFormat = "@"    # keep I31 as text, like the source feed
$ws.Cells.Item(31, 9).Value = "1"    # I31 Antal
$ws.Cells.Item(31, 10).Value = "ex."    # J31 Enhet
$ws.Cells.Item(31, 11).Value = "årsunge"    # K31 Alder-Stadium
$ws.Cells.Item(31, 13).Value = ""    # M31 Aktivitet (now blank but present)
$ws.Cells.Item(31, 17).Value = 561262.631747936    # Q31 Ost
$ws.Cells.Item(31, 18).Value = 6622544.013810508    # R31 Nord

# Row 32
$ws.Cells.Item(32, 1).Value = 111564885    # A32 Id
$ws.Cells.Item(32, 2).Value = 103288    # B32 Taxonsorteringsordning
$ws.Cells.Item(32, 5).Value = 221144    # E32 TaxonId
$ws.Cells.Item(32, 6).Value = "Grönpyrola"    # F32 Artnamn
$ws.Cells.Item(32, 7).Value = "Pyrola chlorantha"    # G32 Vetenskapligt namn
$ws.Cells.Item(32, 8).Value = "Sw."    # H32 Auktor
$ws.Cells.Item(32, 9).NumberFormat = "@"    # keep I32 as text, like the source feed
$ws.Cells.Item(32, 9).Value = "3"    # I32 Antal
$ws.Cells.Item(32, 10).Value = "m²"    # J32 Enhet
$ws.Cells.Item(32, 11).Value = "fullt utvecklade blad"    # K32 Alder-Stadium
$ws.Cells.Item(32, 13).ClearContents()    # M32 Aktivitet (no longer set)
$ws.Cells.Item(32, 17).Value = 561241.7941118333    # Q32 Ost
$ws.Cells.Item(32, 18).Value = 6622674.779475109    # R32 Nord

# Row 33
$ws.Cells.Item(33, 1).Value = 111565024    # A33 Id
$ws.Cells.Item(33, 2).Value = 103288    # B33 Taxonsorteringsordning
$ws.Cells.Item(33, 5).Value = 221144    # E33 TaxonId
$ws.Cells.Item(33, 6).Value = "Grönpyrola"    # F33 Artnamn
$ws.Cells.Item(33, 7).Value = "Pyrola chlorantha"    # G33 Vetenskapligt namn
$ws.Cells.Item(33, 8).Value = "Sw."    # H33 Auktor
$ws.Cells.Item(33, 9).NumberFormat = "@"    # keep I33 as text, like the source feed
$ws.Cells.Item(33, 9).Value = "2"    # I33 Antal
$ws.Cells.Item(33, 10).Value = "m²"    # J33 Enhet
$ws.Cells.Item(33, 11).Value = "fullt utvecklade blad"    # K33 Alder-Stadium
$ws.Cells.Item(33, 17).Value = 561149.6074341368    # Q33 Ost
$ws.Cells.Item(33, 18).Value = 6622721.170183762    # R33 Nord

# Row 34
$ws.Cells.Item(34, 1).Value = 111565017    # A34 Id
$ws.Cells.Item(34, 2).Value = 57578    # B34 Taxonsorteringsordning
$ws.Cells.Item(34, 5).Value = 208250    # E34 TaxonId
$ws.Cells.Item(34, 6).Value = "Åkergroda"    # F34 Artnamn
$ws.Cells.Item(34, 7).Value = "Rana arvalis"    # G34 Vetenskapligt namn
$ws.Cells.Item(34, 8).Value = "Nilsson, 1842"    # H34 Auktor
$ws.Cells.Item(34, 9).NumberFormat = "@"    # keep I34 as text, like the source feed
$ws.Cells.Item(34, 9).Value = "1"    # I34 Antal
$ws.Cells.Item(34, 10).Value = ""    # J34 Enhet
$ws.Cells.Item(34, 11).Value = "årsunge"    # K34 Alder-Stadium
$ws.Cells.Item(34, 13).Value = ""    # M34 Aktivitet (now blank but present)
$ws.Cells.Item(34, 17).Value = 561130.0283522989    # Q34 Ost
$ws.Cells.Item(34, 18).Value = 6622683.03052416    # R34 Nord

# Row 35
$ws.Cells.Item(35, 1).Value = 111565033    # A35 Id
$ws.Cells.Item(35, 2).Value = 103288    # B35 Taxonsorteringsordning
$ws.Cells.Item(35, 5).Value = 221144    # E35 TaxonId
$ws.Cells.Item(35, 6).Value = "Grönpyrola"    # F35 Artnamn
$ws.Cells.Item(35, 7).Value = "Pyrola chlorantha"    # G35 Vetenskapligt namn
$ws.Cells.Item(35, 8).Value = "Sw."    # H35 Auktor
$ws.Cells.Item(35, 9).NumberFormat = "@"    # keep I35 as text, like the source feed
$ws.Cells.Item(35, 9).Value = "2"    # I35 Antal
$ws.Cells.Item(35, 10).Value = "m²"    # J35 Enhet
$ws.Cells.Item(35, 11).Value = "fullt utvecklade blad"    # K35 Alder-Stadium
$ws.Cells.Item(35, 17).Value = 561151.5115810917    # Q35 Ost
$ws.Cells.Item(35, 18).Value = 6622728.260846013    # R35 Nord
